# -----------------------------------------------------------------------
# Fixed update to excel issue
#  1. Rename the "Requested quantity" column header on the two existing
#     sheets to more descriptive names.
#  2. Add a new "PO Forecast" worksheet (after "Monthly Trend") containing
#     the Prophet-style forecast output (ds / PO_Forecast / yhat_lower /
#     yhat_upper).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1. Header renames -------------------------------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. New "PO Forecast" sheet, inserted after the last existing sheet
$lastSheet  = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the sheetPr / print-margin conventions used by the other sheets
$wsForecast.Outline.SummaryRow    = 1
$wsForecast.Outline.SummaryColumn = 1
$wsForecast.PageSetup.LeftMargin   = 54
$wsForecast.PageSetup.RightMargin  = 54
$wsForecast.PageSetup.TopMargin    = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Copy header / date cell formatting from the "Weekly Quantity" sheet so
# the same style indices (bold/centered header, date number format) are
# reused instead of new duplicate styles being created.
$wsWeekly.Range("A1").Copy($wsForecast.Range("A1:D1"))
$wsWeekly.Range("A2").Copy($wsForecast.Range("A2:A21"))

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45011.99999999999, 4,  -9.807782929489219, 16.25320906764547),
    @(45018.99999999999, 4,  -7.975220806966207, 17.13866083452871),
    @(45039.99999999999, 7,  -5.902925996877294, 20.05008063370885),
    @(45046.99999999999, 8,  -5.219092972426244, 21.50839964890155),
    @(45053.99999999999, 9,  -4.517930011416527, 21.92066515664683),
    @(45067.99999999999, 11, -1.832255498419298, 23.64028982069799),
    @(45074.99999999999, 11, -1.647163436808178, 23.99966655099087),
    @(45109.99999999999, 16,  2.361227489599115, 28.91248662049066),
    @(45116.99999999999, 17,  4.143756174505802, 28.91755366826092),
    @(45123.99999999999, 18,  4.317257217880688, 30.53421963434988),
    @(45130.99999999999, 18,  5.385310748365401, 30.72266931649025),
    @(45137.99999999999, 19,  5.615325412112712, 32.93419363834347),
    @(45144.99999999999, 20,  7.123061331883441, 33.40960600684971),
    @(45151.99999999999, 21,  7.158186238467137, 34.74216240997216),
    @(45158.99999999999, 22,  9.122131223200135, 35.6715438298867),
    @(45165.99999999999, 23,  10.22718747551384, 36.05838282156031),
    @(45172.99999999999, 24,  11.53947221532732, 37.18301915164773),
    @(45179.99999999999, 25,  12.04045280439386, 37.80745491508706),
    @(45186.99999999999, 25,  11.6713853148849,  38.63725139101281),
    @(45193.99999999999, 26,  12.76586277555504, 39.87026625769131)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}
